# Refresh market-data-derived Leve profit columns (H-N) across all sheets.
# Values are externally sourced (no formulas in these cells); this mirrors
# the scheduled runner that re-pulls current prices and overwrites them.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 337.2
$ws.Range("J2").Value = 495.16666
$ws.Range("L2").Value = 495.16666
$ws.Range("N2").Value = -721.16666
$ws.Range("H28").Value = 825.6539
$ws.Range("I28").Value = 423.16666
$ws.Range("J28").Value = 1731.25
$ws.Range("K28").Value = 423.16666
$ws.Range("L28").Value = 1731.25
$ws.Range("M28").Value = 61.83334000000002
$ws.Range("N28").Value = -2701.25
$ws.Range("H39").Value = 1271.4
$ws.Range("I39").Value = 281.83334
$ws.Range("K39").Value = 845.5000200000001
$ws.Range("M39").Value = -549.5000200000001
$ws.Range("H62").Value = 8175.2173
$ws.Range("I62").Value = 13878
$ws.Range("K62").Value = 13878
$ws.Range("M62").Value = -13254
$ws.Range("H64").Value = 5129.6665
$ws.Range("I64").Value = 5129.6665
$ws.Range("K64").Value = 5129.6665
$ws.Range("M64").Value = -4881.6665
$ws.Range("H65").Value = 8175.2173
$ws.Range("I65").Value = 13878
$ws.Range("K65").Value = 69390
$ws.Range("M65").Value = -66270
$ws.Range("H67").Value = 5129.6665
$ws.Range("I67").Value = 5129.6665
$ws.Range("K67").Value = 5129.6665
$ws.Range("M67").Value = -4271.6665
$ws.Range("H74").Value = 6303.205
$ws.Range("I74").Value = 3600.625
$ws.Range("J74").Value = 7000.645
$ws.Range("K74").Value = 3600.625
$ws.Range("L74").Value = 7000.645
$ws.Range("M74").Value = -2664.625
$ws.Range("N74").Value = -8872.645
$ws.Range("H77").Value = 6303.205
$ws.Range("I77").Value = 3600.625
$ws.Range("J77").Value = 7000.645
$ws.Range("K77").Value = 18003.125
$ws.Range("L77").Value = 35003.22500000001
$ws.Range("M77").Value = -13323.125
$ws.Range("N77").Value = -44363.22500000001
$ws.Range("H107").Value = 350.8
$ws.Range("I107").Value = 350.8
$ws.Range("K107").Value = 350.8
$ws.Range("M107").Value = 1569.2
$ws.Range("H111").Value = 1158
$ws.Range("I111").Value = 1158
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 3474
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -407
$ws.Range("N111").ClearContents()
$ws.Range("H125").Value = 2476.75
$ws.Range("I125").Value = 2139.8333
$ws.Range("K125").Value = 19258.4997
$ws.Range("M125").Value = -16798.4997
$ws.Range("H132").Value = 2307.9614
$ws.Range("I132").Value = 2200.28
$ws.Range("K132").Value = 6600.84
$ws.Range("M132").Value = -4070.84

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5896.143
$ws.Range("I32").Value = 5896.143
$ws.Range("K32").Value = 5896.143
$ws.Range("M32").Value = -5609.143
$ws.Range("H45").Value = 2615.16
$ws.Range("I45").Value = 2952.5293
$ws.Range("K45").Value = 2952.5293
$ws.Range("M45").Value = -2575.5293
$ws.Range("H61").Value = 3317.8572
$ws.Range("I61").Value = 1448.5454
$ws.Range("J61").Value = 6481.3076
$ws.Range("K61").Value = 1448.5454
$ws.Range("L61").Value = 6481.3076
$ws.Range("M61").Value = -1236.5454
$ws.Range("N61").Value = -6905.3076
$ws.Range("H74").Value = 3110.0435
$ws.Range("I74").Value = 2382.8
$ws.Range("K74").Value = 2382.8
$ws.Range("M74").Value = -1508.8
$ws.Range("H77").Value = 3110.0435
$ws.Range("I77").Value = 2382.8
$ws.Range("K77").Value = 11914
$ws.Range("M77").Value = -7546
$ws.Range("H92").Value = 67958.336
$ws.Range("J92").Value = 67958.336
$ws.Range("L92").Value = 67958.336
$ws.Range("N92").Value = -72950.336
$ws.Range("H110").Value = 662.95
$ws.Range("I110").Value = 644.94116
$ws.Range("K110").Value = 644.94116
$ws.Range("M110").Value = 1400.05884
$ws.Range("H122").Value = 2787.7874
$ws.Range("I122").Value = 2419.0386
$ws.Range("J122").Value = 3244.3333
$ws.Range("K122").Value = 7257.1158
$ws.Range("L122").Value = 9732.999899999999
$ws.Range("M122").Value = -4807.1158
$ws.Range("N122").Value = -14632.9999
$ws.Range("H136").Value = 3317.8572
$ws.Range("I136").Value = 1448.5454
$ws.Range("J136").Value = 6481.3076
$ws.Range("K136").Value = 4345.6362
$ws.Range("L136").Value = 19443.9228
$ws.Range("M136").Value = -1795.6362
$ws.Range("N136").Value = -24543.9228

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1478.6666
$ws.Range("I105").Value = 1456
$ws.Range("K105").Value = 1456
$ws.Range("M105").Value = 291
$ws.Range("H107").Value = 3241.65
$ws.Range("I107").Value = 2696.6155
$ws.Range("K107").Value = 2696.6155
$ws.Range("M107").Value = -776.6154999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3801.4
$ws.Range("I16").Value = 3751.75
$ws.Range("K16").Value = 3751.75
$ws.Range("M16").Value = -3464.75
$ws.Range("H86").Value = 2754
$ws.Range("J86").Value = 2955.3333
$ws.Range("L86").Value = 2955.3333
$ws.Range("N86").Value = -5201.3333
$ws.Range("H89").Value = 2754
$ws.Range("J89").Value = 2955.3333
$ws.Range("L89").Value = 14776.6665
$ws.Range("N89").Value = -26008.6665
$ws.Range("H94").Value = 1104.5358
$ws.Range("J94").Value = 1372.9445
$ws.Range("L94").Value = 1372.9445
$ws.Range("N94").Value = -2274.9445
$ws.Range("H105").Value = 1754.0435
$ws.Range("I105").Value = 1864.7333
$ws.Range("J105").Value = 1546.5
$ws.Range("K105").Value = 1864.7333
$ws.Range("L105").Value = 1546.5
$ws.Range("M105").Value = -117.7333000000001
$ws.Range("N105").Value = -5040.5
$ws.Range("H113").Value = 3801.4
$ws.Range("I113").Value = 3751.75
$ws.Range("K113").Value = 3751.75
$ws.Range("M113").Value = -1581.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 78.666664
$ws.Range("I12").Value = 130.25
$ws.Range("J12").Value = 52.875
$ws.Range("K12").Value = 390.75
$ws.Range("L12").Value = 158.625
$ws.Range("M12").Value = -217.75
$ws.Range("N12").Value = -504.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 79998
$ws.Range("J42").Value = 79998
$ws.Range("L42").Value = 79998
$ws.Range("N42").Value = -80968
$ws.Range("H75").Value = 57999
$ws.Range("J75").Value = 57999
$ws.Range("L75").Value = 57999
$ws.Range("N75").Value = -59747
$ws.Range("H78").Value = 57999
$ws.Range("J78").Value = 57999
$ws.Range("L78").Value = 173997
$ws.Range("N78").Value = -182733
$ws.Range("H80").Value = 6221.6665
$ws.Range("I80").Value = 4434.6665
$ws.Range("K80").Value = 4434.6665
$ws.Range("M80").Value = -3436.6665
$ws.Range("H83").Value = 6221.6665
$ws.Range("I83").Value = 4434.6665
$ws.Range("K83").Value = 22173.3325
$ws.Range("M83").Value = -17181.3325
$ws.Range("H92").Value = 60404.555
$ws.Range("J92").Value = 60404.555
$ws.Range("L92").Value = 60404.555
$ws.Range("N92").Value = -64148.555
$ws.Range("H107").Value = 1472.8
$ws.Range("J107").Value = 2140
$ws.Range("L107").Value = 2140
$ws.Range("N107").Value = -5980
$ws.Range("H113").Value = 3004.6365
$ws.Range("I113").Value = 2208.8
$ws.Range("J113").Value = 3667.8333
$ws.Range("K113").Value = 2208.8
$ws.Range("L113").Value = 3667.8333
$ws.Range("M113").Value = -38.80000000000018
$ws.Range("N113").Value = -8007.8333
$ws.Range("H115").Value = 79998
$ws.Range("J115").Value = 79998
$ws.Range("L115").Value = 79998
$ws.Range("N115").Value = -82348
$ws.Range("H123").Value = 64986.168
$ws.Range("J123").Value = 64986.168
$ws.Range("L123").Value = 64986.168
$ws.Range("N123").Value = -69886.16800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1052.0834
$ws.Range("I55").Value = 243.2
$ws.Range("J55").Value = 1629.8572
$ws.Range("K55").Value = 243.2
$ws.Range("L55").Value = 1629.8572
$ws.Range("M55").Value = -70.19999999999999
$ws.Range("N55").Value = -1975.8572
$ws.Range("H61").Value = 2431.7
$ws.Range("I61").Value = 1530.3636
$ws.Range("J61").Value = 3533.3333
$ws.Range("K61").Value = 1530.3636
$ws.Range("L61").Value = 3533.3333
$ws.Range("M61").Value = -1328.3636
$ws.Range("N61").Value = -3937.3333
$ws.Range("H113").Value = 2431.7
$ws.Range("I113").Value = 1530.3636
$ws.Range("J113").Value = 3533.3333
$ws.Range("K113").Value = 1530.3636
$ws.Range("L113").Value = 3533.3333
$ws.Range("M113").Value = 639.6364000000001
$ws.Range("N113").Value = -7873.3333
$ws.Range("H122").Value = 2938.6365
$ws.Range("I122").Value = 2669.4443
$ws.Range("K122").Value = 8008.3329
$ws.Range("M122").Value = -5558.3329
$ws.Range("H136").Value = 6747.375
$ws.Range("I136").Value = 6121.5
$ws.Range("K136").Value = 18364.5
$ws.Range("M136").Value = -15814.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1282
$ws.Range("I81").Value = 675.5
$ws.Range("K81").Value = 1351
$ws.Range("M81").Value = -290
$ws.Range("H84").Value = 1282
$ws.Range("I84").Value = 675.5
$ws.Range("K84").Value = 6755
$ws.Range("M84").Value = -1451
$ws.Range("H113").Value = 3390.8667
$ws.Range("I113").Value = 1697.125
$ws.Range("K113").Value = 5091.375
$ws.Range("M113").Value = -2921.375
$ws.Range("H136").Value = 7676798.5
$ws.Range("I136").Value = 9335936
$ws.Range("K136").Value = 28007808
$ws.Range("M136").Value = -28005258
